# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff|Handback DateTime"
# timestamps for the 6869bf9d-... (md) row on each sheet, reflecting the
# latest handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 6869bf9d-... file (row 3)
$wsOverview.Range("G3").Value = "2016-09-02 20:56:20"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for
# the 6869bf9d-... file (row 3)
$wsZhCn.Range("H3").Value = "2016-09-02 20:56:15"
$wsZhCn.Range("K3").Value = "2016-09-02 20:56:32"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for
# the 6869bf9d-... file (row 3)
$wsDeDe.Range("H3").Value = "2016-09-02 20:56:20"
$wsDeDe.Range("K3").Value = "2016-09-02 20:56:40"
